# refresh token logic implemented
# Adds two new daily-log rows (07/03 and 07/04) to the Sheet1 task sheet:
#   Row 5: date 45841, hours "8hr 30min", task notes, git repo link (hyperlink)
#   Row 6: date 45842 only (day just started, no other data yet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# keep gridlines visible (matches the source sheet view state; the COM
# round-trip otherwise resets this window-level flag to hidden)
$excel.ActiveWindow.DisplayGridlines = $true

# ---- Row 5 -----------------------------------------------------------
$ws.Range("A5").Value = 45841
$ws.Range("A5").NumberFormat = "dd/mm/yy"
$ws.Range("A5").HorizontalAlignment = -4108

$ws.Range("B5").Value = "8hr 30min"
$ws.Range("B5").HorizontalAlignment = -4108

$taskText = "1)Admin dashboard data fetch `n2)Applying refresh token logic and handling jwt token based on it"
$ws.Range("C5").Value = $taskText
$ws.Range("C5").WrapText = $true

$linkText = "naman-tatvasoft/job-portal (github.com) `nhttps://github.com/naman-tatvasoft/JobApplicationPortal"
$ws.Range("D5").Value = $linkText
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/naman-tatvasoft/JobApplicationPortal", "", "", $linkText)
$ws.Range("D5").WrapText = $true
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 10
$ws.Range("D5").Font.Color = 16711680
$ws.Range("D5").Font.Underline = $false

$ws.Rows.Item(5).RowHeight = 23.85

# ---- Row 6 -----------------------------------------------------------
$ws.Range("A6").Value = 45842
$ws.Range("A6").NumberFormat = "dd/mm/yy"
$ws.Range("A6").HorizontalAlignment = -4108

# ---- selection, matching the final cursor position in the source edit
$ws.Range("C11").Select()
